$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Loan RBI, Variable Instalments: insert a new (blank) column before the
# "Late" column so the repayment schedule gets an extra breakdown column
# (the former N:P "Late"/"Outstanding"/"heading" columns shift right to O:Q).
# A plain column insert takes its width from the column immediately to the
# left (column M), so read that first and re-apply it explicitly - COM's
# Insert() doesn't always stamp a <col> width entry on its own.
$mWidth = $ws.Columns("M").ColumnWidth
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $mWidth

# Make "Repayment schedule" the active sheet/tab and leave the selection on
# G5, matching the saved workbook view state.
$ws.Activate()
$ws.Range("G5").Select()
